# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Fri Nov  3 17:15:06 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Preserve the original style while forcing the numeric-looking
    # string to be stored as text (matches the source data, which
    # keeps prices like "1.02" / "3.20" as literal strings, not numbers).
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "34.904.66"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.828.21"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue $ws.Range("D5") "230.64"
$ws.Range("E5").Value = "  -0.45%  "
Set-TextValue $ws.Range("D6") "0.616"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  +0.13%  "
Set-TextValue $ws.Range("D8") "40.03"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("E9").Value = "  +3.94%  "
Set-TextValue $ws.Range("D10") "0.0683"
$ws.Range("E10").Value = "  +0.36%  "
Set-TextValue $ws.Range("D11") "0.0989"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "2.092.36"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D13") "11.29"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.818.64"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "34.926.85"
$ws.Range("E17").Value = "  +0.28%  "
Set-TextValue $ws.Range("D18") "69.53"
$ws.Range("E18").Value = "  +1.26%  "
Set-TextValue $ws.Range("D20") "239.55"
$ws.Range("E20").Value = "  +1.24%  "
Set-TextValue $ws.Range("D21") "12.15"
$ws.Range("E21").Value = "  +4.01%  "
Set-TextValue $ws.Range("D22") "4.64"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -0.19%  "
Set-TextValue $ws.Range("D25") "173.69"
$ws.Range("E25").Value = "  +0.36%  "
Set-TextValue $ws.Range("D26") "7.74"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  +3.09%  "
Set-TextValue $ws.Range("D28") "17.32"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -5.64%  "
$ws.Range("E30").Value = "  +0.14%  "
Set-TextValue $ws.Range("D31") "0.0550"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  +1.28%  "
Set-TextValue $ws.Range("D33") "3.92"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  +2.48%  "
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("E36").Value = "  +11.36%  "
Set-TextValue $ws.Range("D37") "0.698"
$ws.Range("E37").Value = "  +3.51%  "
Set-TextValue $ws.Range("D38") "92.04"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "1.339.08"
$ws.Range("E39").Value = "  +2.72%  "
Set-TextValue $ws.Range("D40") "1.02"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("E41").Value = "  +1.10%  "
Set-TextValue $ws.Range("D42") "14.49"
$ws.Range("E42").Value = "  -1.95%  "
Set-TextValue $ws.Range("D43") "2.42"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("E44").Value = "  -3.13%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("E46").Value = "  +0.74%  "
Set-TextValue $ws.Range("D47") "0.0522"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").Value = "2.010.21"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  +0.12%  "
Set-TextValue $ws.Range("D50") "0.0668"
$ws.Range("E50").Value = "  +4.05%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D51") "3.20"
$ws.Range("E51").Value = "  +14.01%  "
